$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of rows 2-5 (columns C, D, E, F) before shifting them,
# since the new row2 takes the old row5 numeric data, and rows 3-5 shift up from 2-4.
$origC = @{}
$origD = @{}
$origE = @{}
$origF = @{}
for ($r = 2; $r -le 5; $r++) {
    $origC[$r] = $ws.Cells.Item($r, 3).Value2
    $origD[$r] = $ws.Cells.Item($r, 4).Value2
    $origE[$r] = $ws.Cells.Item($r, 5).Value2
    $origF[$r] = $ws.Cells.Item($r, 6).Value2
}

# Row 2: Technology changes to Biomass_CHP_wood_pellets_DH, numeric data comes from old row 5
$ws.Range("B2").Value = "Biomass_CHP_wood_pellets_DH"
$ws.Range("C2").Value = $origC[5]
$ws.Range("D2").Value = $origD[5]
$ws.Range("E2").Value = $origE[5]
$ws.Range("F2").Value = $origF[5]

# Row 3: numeric data comes from old row 2 (Technology stays Coal PSC)
$ws.Range("C3").Value = $origC[2]
$ws.Range("D3").Value = $origD[2]
$ws.Range("E3").Value = $origE[2]
$ws.Range("F3").Value = $origF[2]

# Row 4: numeric data comes from old row 3 (Technology stays Coal PSC)
$ws.Range("C4").Value = $origC[3]
$ws.Range("D4").Value = $origD[3]
$ws.Range("E4").Value = $origE[3]
$ws.Range("F4").Value = $origF[3]

# Row 5: numeric data comes from old row 4 (Technology stays Coal PSC)
$ws.Range("C5").Value = $origC[4]
$ws.Range("D5").Value = $origD[4]
$ws.Range("E5").Value = $origE[4]
$ws.Range("F5").Value = $origF[4]

# Update Owner column H from "Producer1" to "ProducerNL" for all data rows (2-51)
for ($r = 2; $r -le 51; $r++) {
    if ($ws.Cells.Item($r, 8).Value2 -eq "Producer1") {
        $ws.Cells.Item($r, 8).Value = "ProducerNL"
    }
}
